# Update forecast (预测值, column C) and actual (真实值, column B) values
# for the monthly Natural Gas TTF price dataset on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.7396
$ws.Range("B3").Value = 1.1721
$ws.Range("C3").Value = -1.1538
$ws.Range("B4").Value = 0.2904
$ws.Range("C4").Value = 0.308
$ws.Range("B5").Value = -0.6105
$ws.Range("C5").Value = -0.5764
$ws.Range("B6").Value = -0.6855
$ws.Range("C6").Value = -0.7133
$ws.Range("B7").Value = 0.609
$ws.Range("C7").Value = 0.5826
$ws.Range("B8").Value = -1.9703
$ws.Range("C8").Value = -1.9937
$ws.Range("B9").Value = -3.6198
$ws.Range("C9").Value = -3.5521
$ws.Range("B10").Value = 1.8137
$ws.Range("C10").Value = 1.8072
$ws.Range("B11").Value = 9.5853
$ws.Range("C11").Value = 9.597799999999999
$ws.Range("B12").Value = 7.2314
$ws.Range("C12").Value = 7.201
$ws.Range("B13").Value = 3.8161
$ws.Range("C13").Value = 3.8198
$ws.Range("B14").Value = 3.5875
$ws.Range("C14").Value = 3.5973
$ws.Range("B15").Value = 0.1872
$ws.Range("C15").Value = 0.1728
$ws.Range("B16").Value = -2.399
$ws.Range("C16").Value = -2.3702
$ws.Range("B17").Value = 1.6039
$ws.Range("C17").Value = 1.5822
$ws.Range("C18").Value = -2.1826
$ws.Range("B19").Value = 1.3479
$ws.Range("C19").Value = 1.2711
$ws.Range("B20").Value = 0.2002
$ws.Range("C20").Value = 0.2252
$ws.Range("B21").Value = -2.085
$ws.Range("C21").Value = -2.1584
$ws.Range("B22").Value = -3.7386
$ws.Range("C22").Value = -3.7052
$ws.Range("B23").Value = -5.0153
$ws.Range("C23").Value = -4.9922
$ws.Range("B24").Value = -1.642
$ws.Range("C24").Value = -1.6396
$ws.Range("B25").Value = 3.8974
$ws.Range("B26").Value = 12.7761
$ws.Range("B27").Value = 13.0851
$ws.Range("B28").Value = 2.1504
$ws.Range("B29").Value = -0.2064
$ws.Range("B30").Value = -5.6963
$ws.Range("B31").Value = -2.6631
$ws.Range("B32").Value = -8.606400000000001
$ws.Range("B33").Value = -12.4474
$ws.Range("B34").Value = -23.7046
$ws.Range("B35").Value = -25.5443
$ws.Range("B36").Value = -24.5247
